# Bulk Shift Extension Import.xlsx - add a "Shift" dropdown list on Sheet1!G
# backed by a new helper sheet ("Sheet2") that enumerates the available
# shift names, mirroring the author's change (Bug 562, 571, 567).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the helper sheet and move it right after Sheet1
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Sheet2"
$newSheet.Move($null, $wb.Worksheets.Item("Sheet1"))

# ---------------------------------------------------------------------
# 2. Populate the helper sheet with the list of shift names
# ---------------------------------------------------------------------
$shiftNames = @(
  "Early Shift(06:00-14:30)",
  "Saturday Shift(06:00-11:00)",
  "Morning Shift(07:00-15:30)",
  "General Shift(08:00-16:30)",
  "Sat Shift(08:00-13:00)",
  "Day Shift(09:00-17:30)",
  "Extended Day Shift(10:00-18:30)",
  "Mid Shift(12:00-20:30)",
  "Afternoon Shift(14:00-22:30)",
  "Evening Shift(15:00-00:00)",
  "Late Evening Shift(16:00-01:00)",
  "Late Shift(17:00-02:00)",
  "Night Shift(18:00-03:00)",
  "Extended Night Shift(18:30-03:30)",
  "Late Night Shift(19:30-04:30)",
  "Overnight Shift(21:00-06:00)",
  "Midnight Shift(22:00-07:00)",
  "Weekly Off(00:00-00:00)"
)

$ws2 = $wb.Worksheets.Item("Sheet2")
for ($i = 0; $i -lt $shiftNames.Length; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $shiftNames[$i]
}
$ws2.Columns.Item(1).ColumnWidth = 28.375

# ---------------------------------------------------------------------
# 3. Add the in-cell dropdown (list data validation) on Sheet1 column G,
#    sourced from Sheet2!$A$1:$A$18
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$dvRange = $ws1.Range("G1:G1048576")
$dvRange.Validation.Add(3, 1, 1, 'Sheet2!$A$1:$A$18')
$dvRange.Validation.IgnoreBlank = $true
$dvRange.Validation.InCellDropdown = $true
$dvRange.Validation.ShowInput = $true
$dvRange.Validation.ShowError = $true

# ---------------------------------------------------------------------
# 4. Restore Sheet1's selection, then leave Sheet2 active/selected
#    (last Select()/Activate() wins, matching the saved workbook state)
# ---------------------------------------------------------------------
[void]$ws1.Range("G1").Select()

[void]$ws2.Range("B6").Select()
$ws2.Activate()
